# Natmi following Dr Hou advice
# Update the LR-pair table (Mmp9-Lrp1) on the active sheet:
#  - "ECs" is now included as a Sending cluster (in addition to FAPs, M2, sCs),
#    so the table grows from 12 data rows (rows 2-13) to 16 data rows (rows 2-17).
#  - All numeric columns (E:T) are recomputed accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one full data row: A..T
# Sending cluster, Ligand symbol, Receptor symbol, Target cluster, then 16 numeric metrics
$data = @(
    @("ECs", "Mmp9", "Lrp1", "ECs", 1, 0.3333333333333333, 178.6795653333333, 536.038696, 0.9735959054718317, 0.9735959054718317, 3, 1, 23.63579766666667, 70.907393, 0.06827844587621175, 0.06827844587621175, 4223.234053386614, 38009.10648047952, 0.06647561533705983, 0.06647561533705983),
    @("ECs", "Mmp9", "Lrp1", "FAPs", 1, 0.3333333333333333, 178.6795653333333, 536.038696, 0.9735959054718317, 0.9735959054718317, 3, 1, 181.2883913333334, 543.865174, 0.5237009467675041, 0.523700946767504, 32392.53096341923, 291532.7786707731, 0.5098730974645637, 0.5098730974645637),
    @("ECs", "Mmp9", "Lrp1", "M2", 1, 0.3333333333333333, 178.6795653333333, 536.038696, 0.9735959054718317, 0.9735959054718317, 3, 1, 111.1005463333333, 333.301639, 0.3209442197221123, 0.3209442197221123, 19851.39732713586, 178662.5759442227, 0.3124699782063004, 0.3124699782063004),
    @("ECs", "Mmp9", "Lrp1", "sCs", 1, 0.3333333333333333, 178.6795653333333, 536.038696, 0.9735959054718317, 0.9735959054718317, 3, 1, 30.14303933333333, 90.42911799999999, 0.08707638763417187, 0.08707638763417187, 5385.945165905569, 48473.50649315012, 0.08477721446390778, 0.08477721446390778),
    @("FAPs", "Mmp9", "Lrp1", "ECs", 3, 1, 1.627176, 4.881528, 0.00886621751136806, 0.00886621751136806, 3, 1, 23.63579766666667, 70.907393, 0.06827844587621175, 0.06827844587621175, 38.459602704056, 346.136424336504, 0.000605371552476665, 0.000605371552476665),
    @("FAPs", "Mmp9", "Lrp1", "FAPs", 3, 1, 1.627176, 4.881528, 0.00886621751136806, 0.00886621751136806, 3, 1, 181.2883913333334, 543.865174, 0.5237009467675041, 0.523700946767504, 294.988119456208, 2654.893075105872, 0.004643246504950078, 0.004643246504950076),
    @("FAPs", "Mmp9", "Lrp1", "M2", 3, 1, 1.627176, 4.881528, 0.00886621751136806, 0.00886621751136806, 3, 1, 111.1005463333333, 333.301639, 0.3209442197221123, 0.3209442197221123, 180.780142580488, 1627.021283224392, 0.00284556126107255, 0.00284556126107255),
    @("FAPs", "Mmp9", "Lrp1", "sCs", 3, 1, 1.627176, 4.881528, 0.00886621751136806, 0.00886621751136806, 3, 1, 30.14303933333333, 90.42911799999999, 0.08707638763417187, 0.08707638763417187, 49.048030170256, 441.432271532304, 0.0007720381928687679, 0.0007720381928687679),
    @("M2", "Mmp9", "Lrp1", "ECs", 3, 1, 2.587116666666667, 7.76135, 0.01409677815672808, 0.01409677815672808, 3, 1, 23.63579766666667, 70.907393, 0.06827844587621175, 0.06827844587621175, 61.14856607339444, 550.33709466055, 0.0009625061044031219, 0.000962506104403122),
    @("M2", "Mmp9", "Lrp1", "FAPs", 3, 1, 2.587116666666667, 7.76135, 0.01409677815672808, 0.01409677815672808, 3, 1, 181.2883913333334, 543.865174, 0.5237009467675041, 0.523700946767504, 469.0142186916556, 4221.127968224901, 0.007382496067049965, 0.007382496067049964),
    @("M2", "Mmp9", "Lrp1", "M2", 3, 1, 2.587116666666667, 7.76135, 0.01409677815672808, 0.01409677815672808, 3, 1, 111.1005463333333, 333.301639, 0.3209442197221123, 0.3209442197221123, 287.4300750947389, 2586.87067585265, 0.004524279466106809, 0.004524279466106809),
    @("M2", "Mmp9", "Lrp1", "sCs", 3, 1, 2.587116666666667, 7.76135, 0.01409677815672808, 0.01409677815672808, 3, 1, 30.14303933333333, 90.42911799999999, 0.08707638763417187, 0.08707638763417187, 77.98355944325554, 701.8520349892999, 0.001227496519168181, 0.001227496519168181),
    @("sCs", "Mmp9", "Lrp1", "ECs", 3, 1, 0.631529, 1.894587, 0.003441098860072149, 0.003441098860072149, 3, 1, 23.63579766666667, 70.907393, 0.06827844587621175, 0.06827844587621175, 14.92669166463233, 134.340224981691, 0.0002349528822721302, 0.0002349528822721302),
    @("sCs", "Mmp9", "Lrp1", "FAPs", 3, 1, 0.631529, 1.894587, 0.003441098860072149, 0.003441098860072149, 3, 1, 181.2883913333334, 543.865174, 0.5237009467675041, 0.523700946767504, 114.4888764903487, 1030.399888413138, 0.001802106730940364, 0.001802106730940364),
    @("sCs", "Mmp9", "Lrp1", "M2", 3, 1, 0.631529, 1.894587, 0.003441098860072149, 0.003441098860072149, 3, 1, 111.1005463333333, 333.301639, 0.3209442197221123, 0.3209442197221123, 70.16321692534366, 631.468952328093, 0.001104400788632506, 0.001104400788632506),
    @("sCs", "Mmp9", "Lrp1", "sCs", 3, 1, 0.631529, 1.894587, 0.003441098860072149, 0.003441098860072149, 3, 1, 30.14303933333333, 90.42911799999999, 0.08707638763417187, 0.08707638763417187, 19.03620348714066, 171.325831384266, 0.0002996384582271494, 0.0002996384582271494)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
